$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering for rows 3-20 (Attribute, Type) - rows 1,2,21 remain unchanged
$rows = @{
    3  = @("requested_service_url", "str")
    4  = @("operation_end_time", "datetime")
    5  = @("lifecycle:transition", "str")
    6  = @("complete_service_time", "str")
    7  = @("parameters", "dict")
    8  = @("case:concept:name", "str")
    9  = @("concept:name", "str")
    10 = @("org:resource", "str")
    11 = @("human_workstation_green_button_pressed", "float")
    12 = @("unsatisfied_condition_description", "str")
    13 = @("process_model_id", "str")
    14 = @("case", "str")
    15 = @("SubProcessID", "str")
    16 = @("time:timestamp", "datetime")
    17 = @("identifier:id", "str")
    18 = @("event_id", "str")
    19 = @("current_task", "str")
    20 = @("lifecycle:state", "str")
}

foreach ($r in $rows.Keys) {
    $pair = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}
